# Adding Trending Test Case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 15's Results value changes from PASS to SKIP
$ws.Range("E15").Value = "SKIP"

# New row 16 - TestCase_F15
$ws.Range("A16").Value = "TestCase_F15"
$ws.Range("B16").Value = "OPQA-226"
$ws.Range("C16").Value = "Verify that users should be able to select from a list of suggested topics and check selected topic is presented in users type ahead"
$ws.Range("D16").Value = "Y"
$ws.Range("E16").Value = "PASS"

# Match formatting of the row above (style s="6" for A,B,D,E and s="2" for C)
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("D15:E15").Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# Update selection to match the recorded cursor position
$ws.Range("D9").Select()
